$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 8.453800000000005
$ws.Range("A9").Value = -21.99930000000001
$ws.Range("B9").Value = 5.622900000000008
$ws.Range("D9").Value = -8.273999999999999
$ws.Range("B11").Value = 5.730000000000002
$ws.Range("A18").Value = -22.29420000000001
$ws.Range("A20").Value = -21.31659999999999
$ws.Range("B23").Value = 8.336800000000013
$ws.Range("B24").Value = 4.939100000000004
$ws.Range("B26").Value = 4.950900000000003
$ws.Range("A27").Value = -21.9058
$ws.Range("D27").Value = -8.150100000000002
$ws.Range("D29").Value = -7.409800000000001
$ws.Range("D32").Value = -6.986899999999991
$ws.Range("B34").Value = 9.625100000000002
$ws.Range("A35").Value = -22.0574
$ws.Range("B35").Value = 4.973199999999998
$ws.Range("D37").Value = -7.207000000000001
$ws.Range("D38").Value = -7.777699999999999
$ws.Range("D41").Value = -7.798500000000002
$ws.Range("D45").Value = -7.385599999999996
$ws.Range("B48").Value = 5.746600000000006
$ws.Range("B49").Value = 5.268699999999999
$ws.Range("D51").Value = -8.732899999999995
$ws.Range("B52").Value = 5.776999999999997
$ws.Range("D57").Value = -7.932999999999993
$ws.Range("D64").Value = -7.23819999999999
$ws.Range("B66").Value = 5.241799999999996
$ws.Range("B67").Value = 5.505999999999996
$ws.Range("A69").Value = -21.58569999999999
$ws.Range("A76").Value = -19.58679999999999
$ws.Range("A78").Value = -21.7383
$ws.Range("B78").Value = 5.677599999999997
$ws.Range("B80").Value = 9.326799999999997
$ws.Range("A82").Value = -22.04170000000001
$ws.Range("D82").Value = -9.0764
$ws.Range("A83").Value = -21.56609999999999
$ws.Range("A93").Value = -21.46300000000002
$ws.Range("D93").Value = -6.890499999999992
$ws.Range("B99").Value = 5.968699999999997
$ws.Range("D102").Value = -7.325999999999998
$ws.Range("B104").Value = 9.732700000000001
$ws.Range("D105").Value = -7.576499999999999